# Auto-generated Excel COM-interop script applying cached-value updates
# across 8 worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), per the diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 51458.953
$ws.Range("I15").Value = 51458.953
$ws.Range("K15").Value = 154376.859
$ws.Range("M15").Value = -154207.859

$ws.Range("H17").Value = 301.8
$ws.Range("J17").Value = 312.69446
$ws.Range("L17").Value = 938.08338
$ws.Range("N17").Value = -1274.08338

$ws.Range("H113").Value = 12358790
$ws.Range("J113").Value = 14749.375
$ws.Range("L113").Value = 14749.375
$ws.Range("N113").Value = -21257.375

$ws.Range("H116").Value = 6294.8237
$ws.Range("I116").Value = 2285.7144
$ws.Range("J116").Value = 9101.200000000001
$ws.Range("K116").Value = 2285.7144
$ws.Range("L116").Value = 9101.200000000001
$ws.Range("M116").Value = 1156.2856
$ws.Range("N116").Value = -15985.2

$ws.Range("H129").Value = 2000
$ws.Range("J129").Value = 2000
$ws.Range("L129").Value = 6000
$ws.Range("N129").Value = -16000

$ws.Range("H132").Value = 2120.4707
$ws.Range("I132").Value = 2120.4707
$ws.Range("K132").Value = 6361.4121
$ws.Range("M132").Value = -3831.4121

$ws.Range("H137").Value = 1624.0303
$ws.Range("I137").Value = 1511.48
$ws.Range("J137").Value = 1975.75
$ws.Range("K137").Value = 4534.440000000001
$ws.Range("L137").Value = 5927.25
$ws.Range("M137").Value = -1984.440000000001
$ws.Range("N137").Value = -11027.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 56.666668
$ws.Range("I5").Value = 58
$ws.Range("K5").Value = 58
$ws.Range("M5").Value = 54

$ws.Range("H32").Value = 3381642
$ws.Range("I32").Value = 3679360
$ws.Range("K32").Value = 3679360
$ws.Range("M32").Value = -3679073

$ws.Range("H45").Value = 4502
$ws.Range("I45").Value = 2173
$ws.Range("K45").Value = 2173
$ws.Range("M45").Value = -1796

$ws.Range("H122").Value = 10694.044
$ws.Range("I122").Value = 10948.15
$ws.Range("J122").Value = 9000
$ws.Range("K122").Value = 32844.45
$ws.Range("L122").Value = 27000
$ws.Range("M122").Value = -30394.45
$ws.Range("N122").Value = -31900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 56.666668
$ws.Range("I4").Value = 58
$ws.Range("K4").Value = 58
$ws.Range("M4").Value = 57

$ws.Range("H96").Value = 18971.334

$ws.Range("H99").Value = 3236.2856
$ws.Range("I99").Value = 1672.3077
$ws.Range("K99").Value = 1672.3077
$ws.Range("M99").Value = -174.3077000000001

$ws.Range("H128").Value = 3839.5
$ws.Range("I128").Value = 3839.5
$ws.Range("K128").Value = 11518.5
$ws.Range("M128").Value = -9028.5

$ws.Range("H134").Value = 6032.6113
$ws.Range("I134").Value = 3682.139
$ws.Range("K134").Value = 11046.417
$ws.Range("M134").Value = -8511.417000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1650.1154
$ws.Range("I107").Value = 1388
$ws.Range("J107").Value = 2069.5
$ws.Range("K107").Value = 1388
$ws.Range("L107").Value = 2069.5
$ws.Range("M107").Value = 532
$ws.Range("N107").Value = -5909.5

$ws.Range("H122").Value = 2174.75
$ws.Range("I122").Value = 1500
$ws.Range("K122").Value = 4500
$ws.Range("M122").Value = -2050

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3112303.2
$ws.Range("J4").Value = 3373336.8
$ws.Range("L4").Value = 10120010.4
$ws.Range("N4").Value = -10120234.4

$ws.Range("H5").Value = 1906966.4
$ws.Range("J5").Value = 4482.7144
$ws.Range("L5").Value = 13448.1432
$ws.Range("N5").Value = -13672.1432

$ws.Range("H34").Value = 3863.261
$ws.Range("J34").Value = 7131.6665
$ws.Range("L34").Value = 21394.9995
$ws.Range("N34").Value = -21562.9995

$ws.Range("H62").Value = 6799
$ws.Range("J62").Value = 6799
$ws.Range("L62").Value = 20397
$ws.Range("N62").Value = -21769

$ws.Range("H65").Value = 6799
$ws.Range("J65").Value = 6799
$ws.Range("L65").Value = 61191
$ws.Range("N65").Value = -68055

$ws.Range("H69").Value = 2666.6667
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 2666.6667
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws.Range("H80").Value = 20837624
$ws.Range("J80").Value = 71433144
$ws.Range("L80").Value = 214299432
$ws.Range("N80").Value = -214301304

$ws.Range("H83").Value = 20837624
$ws.Range("J83").Value = 71433144
$ws.Range("L83").Value = 642898296
$ws.Range("N83").Value = -642907656

$ws.Range("H107").Value = 28573398
$ws.Range("J107").Value = 28573398
$ws.Range("L107").Value = 85720194
$ws.Range("N107").Value = -85724034

$ws.Range("H113").Value = 5855.467
$ws.Range("J113").Value = 6794.4165
$ws.Range("L113").Value = 20383.2495
$ws.Range("N113").Value = -24723.2495

$ws.Range("H135").Value = 1906966.4
$ws.Range("J135").Value = 4482.7144
$ws.Range("L135").Value = 40344.4296
$ws.Range("N135").Value = -45414.4296

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 11633
$ws.Range("I35").Value = 2449.5
$ws.Range("K35").Value = 2449.5
$ws.Range("M35").Value = -2151.5

$ws.Range("H113").Value = 5940.4546
$ws.Range("I113").Value = 2972
$ws.Range("J113").Value = 7636.7144
$ws.Range("K113").Value = 2972
$ws.Range("L113").Value = 7636.7144
$ws.Range("M113").Value = -802
$ws.Range("N113").Value = -11976.7144

$ws.Range("H123").Value = 55000
$ws.Range("J123").Value = 55000
$ws.Range("L123").Value = 55000
$ws.Range("N123").Value = -59900

$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530

$ws.Range("H132").Value = 3389.0334
$ws.Range("I132").Value = 3314.3447
$ws.Range("J132").Value = 5555
$ws.Range("K132").Value = 9943.034100000001
$ws.Range("L132").Value = 16665
$ws.Range("M132").Value = -7413.034100000001
$ws.Range("N132").Value = -21725

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 420833.34
$ws.Range("J20").Value = 420833.34
$ws.Range("L20").Value = 420833.34
$ws.Range("N20").Value = -421285.34

$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()

$ws.Range("H61").Value = 4909
$ws.Range("I61").Value = 3399.8462
$ws.Range("K61").Value = 3399.8462
$ws.Range("M61").Value = -3197.8462

$ws.Range("H107").Value = 4998.5
$ws.Range("I107").Value = 4998.5
$ws.Range("K107").Value = 4998.5
$ws.Range("M107").Value = -3078.5

$ws.Range("H113").Value = 4909
$ws.Range("I113").Value = 3399.8462
$ws.Range("K113").Value = 3399.8462
$ws.Range("M113").Value = -1229.8462

$ws.Range("H122").Value = 4476.3184
$ws.Range("J122").Value = 7429.2856
$ws.Range("L122").Value = 22287.8568
$ws.Range("N122").Value = -27187.8568

$ws.Range("H132").Value = 13166639
$ws.Range("I132").Value = 27783596
$ws.Range("J132").Value = 11377.95
$ws.Range("K132").Value = 83350788
$ws.Range("L132").Value = 34133.85000000001
$ws.Range("M132").Value = -83348258
$ws.Range("N132").Value = -39193.85000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

$ws.Range("H81").Value = 8337877.5
$ws.Range("I81").Value = 1752.9
$ws.Range("K81").Value = 3505.8
$ws.Range("M81").Value = -2444.8

$ws.Range("H84").Value = 8337877.5
$ws.Range("I84").Value = 1752.9
$ws.Range("K84").Value = 17529
$ws.Range("M84").Value = -12225

$ws.Range("H126").Value = 2712.625
$ws.Range("I126").Value = 1957.2858
$ws.Range("K126").Value = 5871.857400000001
$ws.Range("M126").Value = -3401.857400000001
